$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing header style (bold, bordered, centered) from A1 onto the
# three new header cells so they match the look of the rest of row 1.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# New header labels for the season-record columns.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Season record values for every player row (2-55).
$ws.Range("AD2:AD55").Value = 97
$ws.Range("AE2:AE55").Value = 65
$ws.Range("AF2:AF55").Value = 0

Write-Host "Season record columns added"
